# Lab 3 Rubric update:
#  - Add "Yes"/"Good" comments into the Grade sheet's Comments (column E) cells
#  - Update window/view state (zoom, scroll position, selection) on both sheets

$wb = $excel.ActiveWorkbook
$wsRubric = $wb.Worksheets.Item("Rubric")
$wsGrade  = $wb.Worksheets.Item("Grade")

# --- Grade sheet: fill in comments column (E) ---
$wsGrade.Range("E6").Value  = "Yes"
$wsGrade.Range("E10").Value = "Yes"
$wsGrade.Range("E11").Value = "Yes"
$wsGrade.Range("E13").Value = "Yes"
$wsGrade.Range("E14").Value = "Yes"
$wsGrade.Range("E16").Value = "Yes"
$wsGrade.Range("E17").Value = "Yes"
$wsGrade.Range("E19").Value = "Good"
$wsGrade.Range("E20").Value = "Yes"
$wsGrade.Range("E21").Value = "Yes"

# --- Grade sheet view: zoom 120%, scroll back to top-left, select B24 ---
$wsGrade.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsGrade.Range("B24").Select()

# --- Rubric sheet view: keep zoom 150%, scroll back to top-left, select D13 ---
$wsRubric.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsRubric.Range("D13").Select()

Write-Output "Rubric updated"
